# Fruta / hortaliza, semanal
# Refresh weekly price data: the source rows were re-pulled and the
# per-row Fecha/Volumen/Precio/Origen figures were shuffled to their
# new reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Registro {
    param([int]$Row, [double]$Fecha, [double]$Volumen, [double]$Precio, [string]$Origen)

    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("M$Row").Value = $Volumen
    $ws.Range("N$Row").Value = $Precio
    $ws.Range("O$Row").Value = $Precio
    $ws.Range("P$Row").Value = $Precio
    $ws.Range("R$Row").Value = $Origen
    $ws.Range("S$Row").Value = $Precio
}

Set-Registro 2  44998 20  2500 "Región de La Araucanía"
Set-Registro 3  44551 120 4500 "Región de O'Higgins"
Set-Registro 4  44215 65  2800 "Región de La Araucanía"
Set-Registro 5  44176 20  3000 "Región de O'Higgins"
Set-Registro 6  44175 40  5000 "Provincia de Curicó"
Set-Registro 7  44574 200 3000 "Región de La Araucanía"
Set-Registro 8  44323 20  3200 "Región de La Araucanía"
Set-Registro 9  44592 5   7500 "Región de La Araucanía"
Set-Registro 10 44616 200 3200 "Región de La Araucanía"
Set-Registro 11 44999 25  2500 "Región de La Araucanía"
Set-Registro 12 44567 80  2400 "Región de La Araucanía"
Set-Registro 13 44214 50  1800 "Región de La Araucanía"
